$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the document binding and action statements with '$' so they
# use the Drools bound-variable syntax ($document instead of document).
$ws.Range("B6").Value = "`$document: Document"
$ws.Range("F7").Value = "`$document.setAutogen(`$1);"
